$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting (incl. the date style on column G) from row 5 down into
# the three new rows, then overwrite with the actual trade data.
$ws.Range("A5:H5").Copy()
$ws.Range("A6:H6").PasteSpecial()
$ws.Range("A7:H7").PasteSpecial()
$ws.Range("A8:H8").PasteSpecial()

# Column A's cached best-fit width shifts slightly (8.85546875 -> 9) now
# that it holds wider numbers (e.g. 10072.64).
$ws.Columns("A:A").ColumnWidth = 8.166666666666666

# Row 6
$ws.Range("A6").Value = 10058.56
$ws.Range("B6").Value = 9992.61
$ws.Range("C6").Value = 305.24
$ws.Range("D6").Value = 307.24
$ws.Range("E6").Value = $false
$ws.Range("F6").Value = 0.66
$ws.Range("G6").Value = 42613.766770833332
$ws.Range("H6").Value = $true

# Row 7
$ws.Range("A7").Value = 10076.67
$ws.Range("B7").Value = 10058.56
$ws.Range("C7").Value = 307.68
$ws.Range("D7").Value = 308.24
$ws.Range("E7").Value = $false
$ws.Range("F7").Value = 0.18
$ws.Range("G7").Value = 42614.674791666665
$ws.Range("H7").Value = $true

# Row 8
$ws.Range("A8").Value = 10072.64
$ws.Range("B8").Value = 10076.67
$ws.Range("C8").Value = 307.95999999999998
$ws.Range("D8").Value = 307.83
$ws.Range("E8").Value = $false
$ws.Range("F8").Value = -0.04
$ws.Range("G8").Value = 42615.751851851855
$ws.Range("H8").Value = $false
